# Improve code for shortcourse MDR improves outcomes
# Insert a new parameter row for "program_prop_treatment_death_shortcoursemdr"
# just above the old row 49 ("program_timeperiod_acf_rounds"), shifting all
# subsequent rows down by one - mirrors the author's edit (a row insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a whole new row at row 49 (pushes old row 49.. down to 50..)
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row with the new parameter + its value
$ws.Range("A49").Value = "program_prop_treatment_death_shortcoursemdr"
$ws.Range("B49").Value = 0.053

# Leave C49/D49/E49 blank, matching the surrounding rows' pattern

# Update the visible selection/scroll position to match where the author ended up
$ws.Activate()
[void]$ws.Range("A35").Select()
$excel.ActiveWindow.ScrollRow = 35
[void]$ws.Range("B50").Select()
